$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "18×43=774" "89×44=3916"
Replace-Text "55×99=5445" "21×26=546"
Replace-Text "71×25=1775" "46×55=2530"
Replace-Text "64×37=2368" "43×80=3440"
Replace-Text "28×83=2324" "62×18=1116"
Replace-Text "17×52=884" "19×32=608"
Replace-Text "22×27=594" "26×99=2574"
Replace-Text "66×91=6006" "58×66=3828"
Replace-Text "33×66=2178" "77×57=4389"
Replace-Text "39×48=1872" "68×89=6052"
Replace-Text "41×53=2173" "50×16=800"
Replace-Text "93×85=7905" "35×67=2345"
Replace-Text "84×73=6132" "61×12=732"
Replace-Text "80×50=4000" "64×81=5184"
Replace-Text "39×96=3744" "74×49=3626"
Replace-Text "54×84=4536" "76×17=1292"
Replace-Text "50×19=950" "56×66=3696"
Replace-Text "85×25=2125" "77×61=4697"
Replace-Text "69×71=4899" "32×17=544"
Replace-Text "34×41=1394" "70×51=3570"
Replace-Text "37×26=962" "88×48=4224"
Replace-Text "85×92=7820" "94×44=4136"
Replace-Text "20×49=980" "98×94=9212"
Replace-Text "93×32=2976" "50×14=700"
Replace-Text "19×67=1273" "15×49=735"
